# "Generate Report for Handback"
#
# The localization-status report tracks handoff/handback state for each
# source file, per target language. The file
# 8c9ec3f7-c2f4-4aca-8b20-1c5c32e828a9.md has now been handed back (in
# sync with en-US) for both zh-cn and de-de, so:
#   - its Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" on the Overview sheet and on each
#     per-language sheet;
#   - each per-language sheet records the new "Latest Handback DateTime"
#     for that file.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# Overview sheet: row 3 is the 8c9ec3f7-...md file.
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# zh-cn sheet: row 3 is the 8c9ec3f7-...md file.
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-03-24 19:03:30"

# de-de sheet: row 3 is the 8c9ec3f7-...md file.
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-03-24 19:03:37"
